$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inputs_for_tollcalib")

# Set MIN_TOLL (column G) to 0.001 (0.1 cents) for all-lane tolling rows 153-319
for ($r = 153; $r -le 319; $r++) {
    $ws.Cells.Item($r, 7).Value = 0.001
}

# Remove the AutoFilter on the data range
$ws.AutoFilterMode = $false

# Update the hidden _FilterDatabase defined name to reflect the new data extent
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Inputs_for_tollcalib!_FilterDatabase") {
        $n.RefersTo = "=Inputs_for_tollcalib!`$B`$1:`$G`$331"
    }
}
